# ============================================================================
# Edit script: "Actualiza base de datos EC y agrega parte 1 de nuevos estado
# de cuenta"
#
# Summary of changes applied to Hoja1:
#  - Header summary values updated (Valor Mora total, Cant. Trabajadores,
#    Cant. Periodos).
#  - The period-detail rows for the existing worker (previously descending
#    2507->2501) are rewritten in ascending order 2501->2507.
#  - A new block of workers (2508 period) is appended below, each occupying
#    one detail row, with the last one keeping the heavier "closing" border
#    style that used to belong to the old last row.
#  - The footer block (signature lines) is moved further down, past the new
#    rows.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# 1) Preserve the footer ("firma") block by copying it down to rows 39-40
#    BEFORE rows 27-28 get overwritten with new worker data below.
# ----------------------------------------------------------------------
$ws.Range("B27:C27").Copy($ws.Range("B39:C39")) | Out-Null
$ws.Range("H27:J27").Copy($ws.Range("H39:J39")) | Out-Null
$ws.Range("B28:C28").Copy($ws.Range("B40:C40")) | Out-Null
$ws.Range("H28:J28").Copy($ws.Range("H40:J40")) | Out-Null

# ----------------------------------------------------------------------
# 2) Preserve the bottom-bordered "closing" row style (currently row 22)
#    by copying it down to row 34, which will become the new last detail
#    row.
# ----------------------------------------------------------------------
$ws.Range("B22:J22").Copy($ws.Range("B34:J34")) | Out-Null

# ----------------------------------------------------------------------
# 3) Stamp the plain detail-row style (currently row 21) across rows
#    22-33 so every new/reordered worker row matches the table look.
# ----------------------------------------------------------------------
for ($r = 22; $r -le 33; $r++) {
    $dest = $ws.Range("B" + $r + ":J" + $r)
    $ws.Range("B21:J21").Copy($dest) | Out-Null
}

# ----------------------------------------------------------------------
# 4) Write the worker/period detail data for rows 16-34.
# ----------------------------------------------------------------------
function Set-RowData {
    param($Row, $TipoDoc, $NumDoc, $Nombre, $Periodo, $ValorMora, $SalarioBasico)
    $ws.Range("B$Row").Value = $TipoDoc
    $ws.Range("C$Row").Value = $NumDoc
    $ws.Range("D$Row").Value = $Nombre
    $ws.Range("E$Row").Value = $Periodo
    $ws.Range("F$Row").Value = $ValorMora
    $ws.Range("G$Row").Value = $SalarioBasico
}

Set-RowData 16 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2501" 43654 1423500
Set-RowData 17 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2502" 56940 1423500
Set-RowData 18 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2503" 56940 1423500
Set-RowData 19 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2504" 56940 1423500
Set-RowData 20 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2505" 56940 1423500
Set-RowData 21 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2506" 56940 1423500
Set-RowData 22 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2507" 56940 1423500
Set-RowData 23 "CC" "73209874" "ROBERT ENRIQUE BORNACHERA DE LA CRUZ" "2508" 56940 1423500
Set-RowData 24 "CC" "45565709" "LUCILA DEL CARMEN MARTINEZ MENDOZA" "2508" 56940 1423500
Set-RowData 25 "CC" "33333428" "OLGA ROQUELINA MARRUGO BARBOSA" "2508" 56940 1423500
Set-RowData 26 "CC" "1002193833" "RAUL ENRIQUE ORTEGA MELENDEZ" "2508" 56940 1423500
Set-RowData 27 "CC" "1047429588" "MARISOL PERTUZ CARREAZO" "2508" 56940 1423500
Set-RowData 28 "CC" "1128060708" "ANGELICA MARGARITA BUELVAS LEON" "2508" 64000 1600000
Set-RowData 29 "CC" "1002476868" "JAIRO ANTONIO TETAY CUADRO" "2508" 56940 1423500
Set-RowData 30 "CC" "1002186795" "ANDRES FELIPE HERNANDEZ MARTINEZ" "2508" 56940 1423500
Set-RowData 31 "CC" "1047477815" "OSMEL DAVID PEÑA GARCIA" "2508" 56940 1423500
Set-RowData 32 "CC" "35896758" "LUZNEY GONZALEZ VALDES" "2508" 56940 1423500
Set-RowData 33 "CC" "1101874137" "ELIZABETH JULIO JULIO" "2508" 56940 1423500
Set-RowData 34 "CC" "1003931691" "HAIDY VERONICA PARRA GONZALEZ" "2508" 56940 1423500

# ----------------------------------------------------------------------
# 5) Refresh the footer text (same wording, new rows 39-40).
# ----------------------------------------------------------------------
$ws.Range("B39").Value = "___________________________________"
$ws.Range("H39").Value = "___________________________________"
$ws.Range("B40").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H40").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# ----------------------------------------------------------------------
# 6) Update the summary header figures.
# ----------------------------------------------------------------------
$ws.Range("E11").Value = 1075634
$ws.Range("C13").Value = 12
$ws.Range("F13").Value = 8

# ----------------------------------------------------------------------
# 7) Widen column D to fit the longest new worker name.
# ----------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 40.5
